$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time per data unit (E2) which cascades into the
# dependent formulas in E3:H3 and E4:H4 via recalculation.
$ws.Range("E2").Value = 20000

$wb.Application.Calculate()
